$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for row 35
$valB35 = 39
$valC35 = $ws.Range("C34").Value()
$valD35 = "39. Excluindo departamentos"
$valE35 = "4:36" + [char]10 + "abordado outra forma de retornar para outra página a partir do controller sem usar o REDIRECT. Invés disso, pode ser usado algum método get do controller, como por exemplo, os metodos de listar todos"

$ws.Range("B35").Value = $valB35
$ws.Range("C35").Value = $valC35
$ws.Range("D35").Value = $valD35
$ws.Range("E35").Value = $valE35

# Re-use the same cell formatting already present in the sheet (style index 5
# for B/C/D, style index 6 -- wrap text -- for E) instead of creating new
# styles: copy format from row 34, which already carries these exact styles.
$ws.Range("B34").Copy()
$ws.Range("B35:D35").PasteSpecial(-4122)
$ws.Range("D34").Copy()
$ws.Range("E35").PasteSpecial(-4122)

# Match the taller row height used for this new note.
$ws.Rows.Item(35).RowHeight = 60

# Update the active selection to E36, as in the target workbook.
[void]$ws.Range("E36").Select()
